$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195; existing rows 195-208 shift down to 196-209.
$ws.Rows.Item(195).Insert()

$ws.Cells.Item(195, 1).Value = 8
$ws.Cells.Item(195, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(195, 3).Value = "Coquimbo"
$ws.Cells.Item(195, 4).Value = 45265
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 6).Value = 100114007
$ws.Cells.Item(195, 7).Value = "Jengibre"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 360
$ws.Cells.Item(195, 11).Value = 22000
$ws.Cells.Item(195, 12).Value = 23000
$ws.Cells.Item(195, 13).Value = 22500
$ws.Cells.Item(195, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(195, 15).Value = "Perú"
$ws.Cells.Item(195, 16).Value = 1731
$ws.Cells.Item(195, 17).Value = 13
$ws.Cells.Item(195, 18).Value = "Hortaliza"
